$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Formula = "1.438467860221863"
$ws.Range("B1").Formula = "3.517448663711548"
$ws.Range("C1").Formula = "5.267248630523682"
$ws.Range("D1").Formula = "1.723904609680176"
$ws.Range("E1").Formula = "0.9659792184829712"
